$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from SCD0201 to SCD0011
$ws.Name = "SCD0011"

# Update Test Case ID in B2 from "DGS-216" to "SCD0011-032"
$ws.Range("B2").Value = "SCD0011-032"

# Widen column B to fit new content (closest achievable value to target stored width 12.42578125)
$ws.Columns("B").ColumnWidth = 11.6

# Update sheet view / selection: remove frozen/scrolled topLeftCell, select B3 instead of P2
$ws.Range("B3").Select()
